$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.959.42"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "'2.943.68"
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'552.86"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "'133.09"
$ws.Range("E6").Value = "  +10.53%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +5.71%  "

$ws.Range("D9").Value = "'2.935.40"
$ws.Range("E9").Value = "  +3.00%  "

$ws.Range("E10").Value = "  +3.72%  "

$ws.Range("D11").Value = "'4.81"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = "  +5.61%  "

$ws.Range("E13").Value = "  +5.83%  "

$ws.Range("D14").Value = "'32.79"
$ws.Range("E14").Value = "  +6.16%  "

$ws.Range("E15").Value = "  +3.94%  "

$ws.Range("D16").Value = "'3.428.95"
$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").Value = "'6.89"
$ws.Range("E17").Value = "  +10.97%  "

$ws.Range("D18").Value = "'2.940.49"
$ws.Range("E18").Value = "  +3.01%  "

$ws.Range("D19").Value = "'57.950.23"
$ws.Range("E19").Value = "  +1.40%  "

$ws.Range("D20").Value = "'416.81"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").Value = "'13.31"
$ws.Range("E21").Value = "  +5.88%  "

$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  +8.54%  "

$ws.Range("D23").Value = "'13.43"
$ws.Range("E23").Value = "  +9.05%  "

$ws.Range("D24").Value = "'6.99"
$ws.Range("E24").Value = "  +5.22%  "

$ws.Range("D25").Value = "'78.92"
$ws.Range("E25").Value = "  +4.47%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("E28").Value = "  +3.59%  "

$ws.Range("D29").Value = "'2.02"
$ws.Range("E29").Value = "  +7.55%  "

$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = "  +6.52%  "

$ws.Range("E31").Value = "  +4.64%  "

$ws.Range("D32").Value = "'5.93"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").Value = "'0.0971"
$ws.Range("E33").Value = "  +6.18%  "

$ws.Range("D34").Value = "'5.68"
$ws.Range("E34").Value = "  +7.54%  "

$ws.Range("D35").Value = "'0.946"
$ws.Range("E35").Value = "  +8.67%  "

$ws.Range("E36").Value = "  +5.63%  "

$ws.Range("D37").Value = "'0.0₃0702"
$ws.Range("E37").Value = "  +15.53%  "

$ws.Range("D38").Value = "'48.23"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").Value = "'8.71"
$ws.Range("E39").Value = "  +7.43%  "

$ws.Range("D40").Value = "'2.65"
$ws.Range("E40").Value = "  +14.23%  "

$ws.Range("D41").Value = "'379.89"
$ws.Range("E41").Value = "  +10.67%  "

$ws.Range("E42").Value = "  +5.43%  "

$ws.Range("D43").Value = "'0.0347"
$ws.Range("E43").Value = "  +3.85%  "

$ws.Range("D44").Value = "'2.696.65"
$ws.Range("E44").Value = "  +4.84%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").Value = "'123.90"
$ws.Range("E46").Value = "  +6.27%  "

$ws.Range("D47").Value = "'0.236"
$ws.Range("E47").Value = "  +5.41%  "

$ws.Range("D48").Value = "'1.96"
$ws.Range("E48").Value = "  +3.81%  "

$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").Value = "'22.89"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("D51").Value = "'1.99"
$ws.Range("E51").Value = "  +5.07%  "
